# Insert a new weekly price-report row at row 13, pushing the existing
# rows 13..101 down to 14..102, and fill the new row with the latest
# market observation (Poroto verde / Magnum, "malla 25 kilos" ex Peru).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 13 (shifts 13:101 -> 14:102).
$ws.Rows(13).Insert()

# Populate the newly inserted row 13 with the new record.
$ws.Cells.Item(13, 1).Value  = 7
$ws.Cells.Item(13, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(13, 3).Value  = "Ñuble"
$ws.Cells.Item(13, 4).Value  = 44859
$ws.Cells.Item(13, 5).Value  = 16
$ws.Cells.Item(13, 6).Value  = 100112031
$ws.Cells.Item(13, 7).Value  = "Poroto verde"
$ws.Cells.Item(13, 8).Value  = "Magnum"
$ws.Cells.Item(13, 9).Value  = "Primera"
$ws.Cells.Item(13, 10).Value = 60
$ws.Cells.Item(13, 11).Value = 32000
$ws.Cells.Item(13, 12).Value = 33000
$ws.Cells.Item(13, 13).Value = 32500
$ws.Cells.Item(13, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(13, 15).Value = "Perú"
$ws.Cells.Item(13, 16).Value = 1300
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = "Hortaliza"
